$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'26.571.66"
$ws.Cells.Item(2, 5).Value = "'  +0.80%  "
$ws.Cells.Item(3, 4).Value = "'1.731.92"
$ws.Cells.Item(3, 5).Value = "'  +0.46%  "
$ws.Cells.Item(4, 4).Value = "'0.9985"
$ws.Cells.Item(4, 5).Value = "'  -0.04%  "
$ws.Cells.Item(5, 4).Value = "'245.22"
$ws.Cells.Item(5, 5).Value = "'  +1.36%  "
$ws.Cells.Item(6, 4).Value = "'0.9991"
$ws.Cells.Item(6, 5).Value = "'  -0.04%  "
$ws.Cells.Item(7, 4).Value = "'0.4816"
$ws.Cells.Item(7, 5).Value = "'  +1.69%  "
$ws.Cells.Item(8, 4).Value = "'0.2686"
$ws.Cells.Item(8, 5).Value = "'  +1.88%  "
$ws.Cells.Item(9, 4).Value = "'0.06193"
$ws.Cells.Item(9, 5).Value = "'  -0.19%  "
$ws.Cells.Item(10, 4).Value = "'1.727.57"
$ws.Cells.Item(10, 5).Value = "'  +0.27%  "
$ws.Cells.Item(11, 4).Value = "'0.07192"
$ws.Cells.Item(11, 5).Value = "'  +1.75%  "
$ws.Cells.Item(12, 4).Value = "'15.61"
$ws.Cells.Item(12, 5).Value = "'  +0.70%  "
$ws.Cells.Item(13, 4).Value = "'0.6120"
$ws.Cells.Item(13, 5).Value = "'  +2.41%  "
$ws.Cells.Item(14, 4).Value = "'4.540"
$ws.Cells.Item(14, 5).Value = "'  +2.54%  "
$ws.Cells.Item(15, 4).Value = "'77.42"
$ws.Cells.Item(15, 5).Value = "'  +1.31%  "
$ws.Cells.Item(16, 4).Value = "'0.9993"
$ws.Cells.Item(16, 5).Value = "'  -0.01%  "
$ws.Cells.Item(17, 4).Value = "'26.565.21"
$ws.Cells.Item(17, 5).Value = "'  +0.73%  "
$ws.Cells.Item(18, 4).Value = "'0.9985"
$ws.Cells.Item(18, 5).Value = "'  -0.10%  "
$ws.Cells.Item(19, 4).Value = "'0.000006972"
$ws.Cells.Item(19, 5).Value = "'  +2.03%  "
$ws.Cells.Item(20, 4).Value = "'11.57"
$ws.Cells.Item(20, 5).Value = "'  +0.02%  "
$ws.Cells.Item(21, 4).Value = "'1.951.33"
$ws.Cells.Item(21, 5).Value = "'  +0.58%  "
$ws.Cells.Item(22, 4).Value = "'4.530"
$ws.Cells.Item(22, 5).Value = "'  +0.17%  "
$ws.Cells.Item(23, 4).Value = "'8.825"
$ws.Cells.Item(23, 5).Value = "'  +0.79%  "
$ws.Cells.Item(24, 4).Value = "'5.265"
$ws.Cells.Item(24, 5).Value = "'  +0.07%  "
$ws.Cells.Item(25, 4).Value = "'137.15"
$ws.Cells.Item(25, 5).Value = "'  +1.39%  "
$ws.Cells.Item(26, 4).Value = "'15.41"
$ws.Cells.Item(26, 5).Value = "'  +1.04%  "
$ws.Cells.Item(27, 4).Value = "'1.786"
$ws.Cells.Item(27, 5).Value = "'  +0.85%  "
$ws.Cells.Item(28, 4).Value = "'1.409"
$ws.Cells.Item(28, 5).Value = "'  +0.61%  "
$ws.Cells.Item(29, 4).Value = "'107.72"
$ws.Cells.Item(29, 5).Value = "'  +0.68%  "
$ws.Cells.Item(30, 4).Value = "'3.991"
$ws.Cells.Item(30, 5).Value = "'  +0.78%  "
$ws.Cells.Item(31, 4).Value = "'0.08039"
$ws.Cells.Item(31, 5).Value = "'  +2.95%  "
$ws.Cells.Item(32, 4).Value = "'3.706"
$ws.Cells.Item(32, 5).Value = "'  +0.45%  "
$ws.Cells.Item(33, 4).Value = "'0.04523"
$ws.Cells.Item(33, 5).Value = "'  +0.47%  "
$ws.Cells.Item(34, 4).Value = "'2.617"
$ws.Cells.Item(34, 5).Value = "'  +0.05%  "
$ws.Cells.Item(35, 4).Value = "'1.011"
$ws.Cells.Item(35, 5).Value = "'  +2.92%  "
$ws.Cells.Item(36, 4).Value = "'0.6279"
$ws.Cells.Item(36, 5).Value = "'  +0.74%  "
$ws.Cells.Item(37, 4).Value = "'2.095"
$ws.Cells.Item(37, 5).Value = "'  +8.17%  "
$ws.Cells.Item(38, 4).Value = "'0.9103"
$ws.Cells.Item(38, 5).Value = "'  -3.77%  "
$ws.Cells.Item(39, 4).Value = "'2.391"
$ws.Cells.Item(39, 5).Value = "'  -2.64%  "
$ws.Cells.Item(40, 4).Value = "'1.001"
$ws.Cells.Item(40, 5).Value = "'  +0.01%  "
$ws.Cells.Item(41, 4).Value = "'0.01506"
$ws.Cells.Item(41, 5).Value = "'  +1.27%  "
$ws.Cells.Item(42, 4).Value = "'102.51"
$ws.Cells.Item(42, 5).Value = "'  -10.20%  "
$ws.Cells.Item(43, 4).Value = "'5.542"
$ws.Cells.Item(43, 5).Value = "'  -2.66%  "
$ws.Cells.Item(44, 4).Value = "'0.3891"
$ws.Cells.Item(44, 5).Value = "'  +1.52%  "
$ws.Cells.Item(45, 4).Value = "'7.035"
$ws.Cells.Item(45, 5).Value = "'  +10.45%  "
$ws.Cells.Item(46, 4).Value = "'0.1183"
$ws.Cells.Item(46, 5).Value = "'  -0.28%  "
$ws.Cells.Item(47, 4).Value = "'0.05375"
$ws.Cells.Item(47, 5).Value = "'  +1.92%  "
$ws.Cells.Item(48, 4).Value = "'30.77"
$ws.Cells.Item(48, 5).Value = "'  +1.13%  "
$ws.Cells.Item(49, 4).Value = "'7.844"
$ws.Cells.Item(49, 5).Value = "'  -0.31%  "
$ws.Cells.Item(50, 4).Value = "'1.253"
$ws.Cells.Item(50, 5).Value = "'  +2.69%  "
$ws.Cells.Item(51, 4).Value = "'0.3412"
$ws.Cells.Item(51, 5).Value = "'  +0.77%  "
